$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 122
$ws.Range("F6").Value = 2959
$ws.Range("F8").Value = 1969
$ws.Range("F11").Value = 802
$ws.Range("F12").Value = 940
$ws.Range("F13").Value = 195
$ws.Range("F14").Value = 410
$ws.Range("F15").Value = 1133
$ws.Range("F17").Value = 63
$ws.Range("F19").Value = 7112
$ws.Range("F20").Value = 266
$ws.Range("F21").Value = 1809
$ws.Range("F23").Value = 186
$ws.Range("F25").Value = 372
$ws.Range("F26").Value = 316
$ws.Range("F27").Value = 76
$ws.Range("F28").Value = 1117
$ws.Range("F29").Value = 940
$ws.Range("F31").Value = 115
$ws.Range("F33").Value = 934
$ws.Range("F34").Value = 1916
$ws.Range("F35").Value = 467
$ws.Range("F36").Value = 6
$ws.Range("F37").Value = 152
$ws.Range("F38").Value = 246
$ws.Range("F39").Value = 33
$ws.Range("F41").Value = 266

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 13

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 122
$ws.Range("F6").Value = 13
$ws.Range("F9").Value = 2959
$ws.Range("F11").Value = 1969
$ws.Range("F14").Value = 802
$ws.Range("F16").Value = 940
$ws.Range("F17").Value = 195
$ws.Range("F18").Value = 410
$ws.Range("F19").Value = 1133
$ws.Range("F21").Value = 63
$ws.Range("F23").Value = 7112
$ws.Range("F24").Value = 266
$ws.Range("F25").Value = 1810
$ws.Range("F28").Value = 186
$ws.Range("F30").Value = 372
$ws.Range("F31").Value = 316
$ws.Range("F32").Value = 76
$ws.Range("F33").Value = 1117
$ws.Range("F34").Value = 940
$ws.Range("F35").Value = 63
$ws.Range("F36").Value = 115
$ws.Range("F37").Value = 934
$ws.Range("F38").Value = 1916
$ws.Range("F39").Value = 468
$ws.Range("F40").Value = 6
$ws.Range("F41").Value = 152
$ws.Range("F42").Value = 246
$ws.Range("F43").Value = 33
$ws.Range("F45").Value = 266
